# Remove the "contextual spacing" direct-formatting override
# (w:contextualSpacing, i.e. the "Don't add space between paragraphs of
# the same style" paragraph setting) from every paragraph in the
# document body as well as from the paragraph(s) inside the comments.
#
# This mirrors Word's own behaviour: selecting the paragraphs, opening
# Format > Paragraph, and leaving "Don't add space between paragraphs
# of the same style" unset clears the explicit <w:contextualSpacing/>
# element that a previous pass had stamped onto every paragraph.

$d = $word.ActiveDocument

function Clear-ContextualSpacing($paragraphs) {
    foreach ($p in $paragraphs) {
        try {
            $p.Format.ContextualSpacing = $false
        } catch {
            # Older/limited hosts may not surface this member; ignore and
            # continue so the rest of the run still applies.
        }
    }
}

# Body paragraphs.
Clear-ContextualSpacing $d.Paragraphs

# Paragraphs living inside comments (the comment text itself carries the
# same paragraph-level override in this document).
foreach ($cmt in $d.Comments) {
    try {
        Clear-ContextualSpacing $cmt.Range.Paragraphs
    } catch {
    }
}

Write-Output "contextualSpacing cleared"
